# Adds a "2022-Q4" quarterly sheet to the "300638-广和通" fund-holdings
# workbook and records its summary row in the "总计" (totals) sheet.
#
# Before: 总计, 2022-Q3, 2022-Q2, 2022-Q1, 2021-Q4, 2021-Q3, 2021-Q2, 2021-Q1, 2020-Q4
# After:  总计, 2022-Q4, 2022-Q3, 2022-Q2, 2022-Q1, 2021-Q4, 2021-Q3, 2021-Q2, 2021-Q1, 2020-Q4

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update the "总计" sheet: insert a new row for 2022-Q4 right after the
#    header row, pushing every existing quarter down by one row.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q4"
$total.Cells.Item(2, 3).Value = 4
$total.Cells.Item(2, 4).Value = 0.02

# Renumber the "持有数量(只)" index column (col A) for the rows that shifted
# down, so it stays a plain 0-based running sequence (0,1,2,...).
For ($r = 3; $r -le 10; $r++) {
    $total.Cells.Item($r, 1).Value = $r - 2
}

# ---------------------------------------------------------------------------
# 2. Insert the new "2022-Q4" worksheet right after "总计" (so it lands
#    before "2022-Q3", matching the tab order in the diff) and fill it with
#    the quarter's fund-holdings data.
# ---------------------------------------------------------------------------
$q4 = $wb.Worksheets.Add($null, $total)
$q4.Name = "2022-Q4"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
For ($c = 0; $c -lt $headers.Length; $c++) {
    $q4.Cells.Item(1, $c + 2).Value = $headers[$c]
}

$rows = @(
    @("013733", "红塔红土盛丰混合A",         "0.40", "61.27", "2.91", "0.0116", 9),
    @("005502", "华泰紫金智能量化股票A",       "0.24", "94.33", "1.16", "0.0028", 6),
    @("013734", "红塔红土盛丰混合C",         "0.09", "61.27", "2.91", "0.0026", 9),
    @("014629", "华泰紫金智能量化股票C",       "0.00", "94.33", "1.16", 0,        6)
)

For ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $row = $rows[$i]
    $q4.Cells.Item($r, 1).Value = $i

    # Fund code / name stay plain text.
    $q4.Cells.Item($r, 2).Value = "'" + $row[0]
    $q4.Cells.Item($r, 3).Value = $row[1]

    # 基金规模 / 股票总仓位 / 仓位占比 are stored as text in the source data
    # (leading apostrophe keeps e.g. "0.40" / "0.00" from being normalised
    # to a number and losing the trailing zero).
    $q4.Cells.Item($r, 4).Value = "'" + $row[2]
    $q4.Cells.Item($r, 5).Value = "'" + $row[3]
    $q4.Cells.Item($r, 6).Value = "'" + $row[4]

    # 持有市值(亿元): numeric 0 for the last row, text for the others.
    if ($row[5] -eq 0) {
        $q4.Cells.Item($r, 7).Value = 0
    } else {
        $q4.Cells.Item($r, 7).Value = "'" + $row[5]
    }

    # 仓位排名 is numeric.
    $q4.Cells.Item($r, 8).Value = $row[6]
}
